# Update "想去人数" (interested-count) figures scraped at commit 456a3b4.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F9").Value  = 7300
$ws1.Range("F37").Value = 136
$ws1.Range("F45").Value = 511
$ws1.Range("F46").Value = 241
$ws1.Range("F48").Value = 735
$ws1.Range("F50").Value = 98

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F29").Value = 19

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F8").Value  = 2883
$ws3.Range("F9").Value  = 1123
$ws3.Range("F10").Value = 1096
$ws3.Range("F13").Value = 1812
$ws3.Range("F14").Value = 8114

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value  = 2883
$ws4.Range("F7").Value  = 7300
$ws4.Range("F8").Value  = 1123
$ws4.Range("F36").Value = 136
$ws4.Range("F45").Value = 241
$ws4.Range("F47").Value = 19
